$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# The insert copies formatting from row 12 into A13; clear it since target has no A13 cell
$ws.Range("A13").Clear()

# Give B13/C13 the correct number formats/styles (copy from B14/C14, which already carry styles 2/3)
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new docente row
$ws.Range("B13").Value = '9146830 - Danúbia Caporusso Bargos'
$ws.Range("C13").Value = '9146830 - Danúbia Caporusso Bargos'

# Update text content for the shifted rows (14-24)
$ws.Range("B14").Value = 'Elementos de teoria e história do planejamento urbano. Teoria e prática do planejamento ambiental; Planejamento ambiental como indutor de desenvolvimento sustentável; Aplicações da teoria do planejamento a problemas ambientais e urbanos; Legislação e política ambiental urbana.'
$ws.Range("C14").Value = 'Elementos de teoria e história do planejamento urbano. Teoria e prática do planejamento ambiental; Planejamento ambiental como indutor de desenvolvimento sustentável; Aplicações da teoria do planejamento a problemas ambientais e urbanos; Legislação e política ambiental urbana.'

$ws.Range("B15").Value = 'Theory elements and history of urban planning. Theory and practice of environmental planning; environmental planning as an inducer of sustainable development; environmental theory applied to urban problems; legislation and urban environmental policy.'
$ws.Range("C15").Value = 'Theory elements and history of urban planning. Theory and practice of environmental planning; environmental planning as an inducer of sustainable development; environmental theory applied to urban problems; legislation and urban environmental policy.'

$ws.Range("B16").Value = 'Introdução ao planejamento e gestão ambiental. Origens da teoria e prática do planejamento. Natureza do planejamento e suas relações com a geografia, política, economia, sociedade, cultura e meio ambiente. Análises, estudos e proposições relativas às diversas formas de crescimento e expansão urbanas; Elementos para estruturação ambiental da cidade; Etapas, estruturas e instrumentos do planejamento ambiental; Indicadores ambientais e planejamento; Participação pública no planejamento ambiental; Política Nacional do Meio Ambiente (Lei n°6938/1981); Sistema Nacional de Unidades de Conservação (Lei n°9985/2000); Estatuto da Cidade (Lei n°10.257/2001); Zoneamento Ambiental; EIA e EIV como instrumentos inovadores; Novos conceitos e princípios de planos diretores urbano-ambientais;'
$ws.Range("C16").Value = 'Introdução ao planejamento e gestão ambiental. Origens da teoria e prática do planejamento. Natureza do planejamento e suas relações com a geografia, política, economia, sociedade, cultura e meio ambiente. Análises, estudos e proposições relativas às diversas formas de crescimento e expansão urbanas; Elementos para estruturação ambiental da cidade; Etapas, estruturas e instrumentos do planejamento ambiental; Indicadores ambientais e planejamento; Participação pública no planejamento ambiental; Política Nacional do Meio Ambiente (Lei n°6938/1981); Sistema Nacional de Unidades de Conservação (Lei n°9985/2000); Estatuto da Cidade (Lei n°10.257/2001); Zoneamento Ambiental; EIA e EIV como instrumentos inovadores; Novos conceitos e princípios de planos diretores urbano-ambientais;'

$ws.Range("B17").Value = 'Environmental planning and management introduction. Planning theory and practice origins. Nature of planning and its relations with geography, politics, economy, society, culture and environment. Analyzes, studies and propositions related to different types of urban growth and expansion; Elements for city environmental structuring; Stages, structures and instruments of environmental planning; Environmental indicators and planning; Public participation in environmental planning; National Policy of the Environment; National System of Conservation Units; City Statute; Environmental Zoning; EIA and EIV as innovative instruments; New concepts and principles of urban-environmental master plans;.'
$ws.Range("C17").Value = 'Environmental planning and management introduction. Planning theory and practice origins. Nature of planning and its relations with geography, politics, economy, society, culture and environment. Analyzes, studies and propositions related to different types of urban growth and expansion; Elements for city environmental structuring; Stages, structures and instruments of environmental planning; Environmental indicators and planning; Public participation in environmental planning; National Policy of the Environment; National System of Conservation Units; City Statute; Environmental Zoning; EIA and EIV as innovative instruments; New concepts and principles of urban-environmental master plans;.'

$ws.Range("B19").Value = 'Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C19").Value = 'Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'

$ws.Range("B20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'

$ws.Range("B21").Value = 'Provas e/ou exercícios dirigidos.'
$ws.Range("C21").Value = 'Provas e/ou exercícios dirigidos.'

$ws.Range("B22").Value = 'Bibliografia básica:AGRA FILHO, S,S. Planejamento e Gestão Ambiental no Brasil. Os Instrumentos da Política Nacional do Meio Ambiente, Rio de Janeiro, Elsevier, 2014FRANCO, M.A.R., Planejamento ambiental para a cidade sustentável, Ed. Annablume, 2000DEAK, C., SHIFFER, S.T.R., O processo de urbanização no Brasil, EDUSP, 1999IBGE, Instituto Brasileiro de Geografia e Estatística. Indicadores de Desenvolvimento Sustentável. Rio de Janeiro, IBGE, 2012.MOTA, S., Urbanização e meio ambiente, ABES Associação Brasileira de Engenharia Sanitária, 1999MENEZES, C.L., Desenvolvimento urbano e meio ambiente, Papirus, 1996PHILLIPI, Jr.A; MALHEIROS, T.F. Indicadores de Sustentabilidade e Gestão Ambiental. Editora Manole, 2012.SANTOS, M. A Urbanização Brasileira. 3 ed. São Paulo: HUCITEC, 1993. 155pSANTOS, R.F., Planejamento ambiental: teoria e prática, Editora Oficina de textos, 2004SECCHI, L. Análise de Políticas Públicas. Diagnóstico de Problemas, Recomendações de Soluções., São Paulo, Cengage Learning, 2016SOUZA, M.L. Mudar a Cidade: Uma introdução crítica ao planejamento e à gestão urbanos. Rio de Janeiro, Bertrand Brasil, 2003.VILLAÇA, F. Uma contribuição para a história do planejamento urbano no Brasil. In: DEAK, C; SCHIFFER, S.R (org) O processo de urbanização no Brasil. São Paulo, EDUSP, 1999.Bibliografia complementar:ALLEN, A., YOU, N., Sustainable urbanization – bridging the green and brown agendas, DPU, University College London, 2002ACSELRAD, H., Conflitos ambientais no Brasil, Fundação Henrich Boll, 2004BARDET, G., O urbanismo, Papirus, 1990BUARQUE, S.C., LIMA, R.R.A.; Manual de estratégia de desenvolvimento para aglomerações urbanas, Brasília, IPEA, 2005MENEGAT, R; ALMEIDA, G. Desenvolvimento Sustentável e Gestão Ambiental nas Cidades. Porto Alegre, Editora UFRGS, 2004.'
$ws.Range("C22").Value = 'Bibliografia básica:AGRA FILHO, S,S. Planejamento e Gestão Ambiental no Brasil. Os Instrumentos da Política Nacional do Meio Ambiente, Rio de Janeiro, Elsevier, 2014FRANCO, M.A.R., Planejamento ambiental para a cidade sustentável, Ed. Annablume, 2000DEAK, C., SHIFFER, S.T.R., O processo de urbanização no Brasil, EDUSP, 1999IBGE, Instituto Brasileiro de Geografia e Estatística. Indicadores de Desenvolvimento Sustentável. Rio de Janeiro, IBGE, 2012.MOTA, S., Urbanização e meio ambiente, ABES Associação Brasileira de Engenharia Sanitária, 1999MENEZES, C.L., Desenvolvimento urbano e meio ambiente, Papirus, 1996PHILLIPI, Jr.A; MALHEIROS, T.F. Indicadores de Sustentabilidade e Gestão Ambiental. Editora Manole, 2012.SANTOS, M. A Urbanização Brasileira. 3 ed. São Paulo: HUCITEC, 1993. 155pSANTOS, R.F., Planejamento ambiental: teoria e prática, Editora Oficina de textos, 2004SECCHI, L. Análise de Políticas Públicas. Diagnóstico de Problemas, Recomendações de Soluções., São Paulo, Cengage Learning, 2016SOUZA, M.L. Mudar a Cidade: Uma introdução crítica ao planejamento e à gestão urbanos. Rio de Janeiro, Bertrand Brasil, 2003.VILLAÇA, F. Uma contribuição para a história do planejamento urbano no Brasil. In: DEAK, C; SCHIFFER, S.R (org) O processo de urbanização no Brasil. São Paulo, EDUSP, 1999.Bibliografia complementar:ALLEN, A., YOU, N., Sustainable urbanization – bridging the green and brown agendas, DPU, University College London, 2002ACSELRAD, H., Conflitos ambientais no Brasil, Fundação Henrich Boll, 2004BARDET, G., O urbanismo, Papirus, 1990BUARQUE, S.C., LIMA, R.R.A.; Manual de estratégia de desenvolvimento para aglomerações urbanas, Brasília, IPEA, 2005MENEGAT, R; ALMEIDA, G. Desenvolvimento Sustentável e Gestão Ambiental nas Cidades. Porto Alegre, Editora UFRGS, 2004.'

$ws.Range("B24").Value = 'LOB1235 -  Avaliação de Impactos Ambientais  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOB1235 -  Avaliação de Impactos Ambientais  (Requisito fraco)
'

# Fix up row heights to match the final layout
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30

